$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87
$ws.Cells.Item($row, 1).Value = 221
$ws.Cells.Item($row, 2).Value = 138
$ws.Cells.Item($row, 3).Value = 76
$ws.Cells.Item($row, 4).Value = 4
$ws.Cells.Item($row, 5).Value = 3
$ws.Cells.Item($row, 6).Value = 86
$ws.Cells.Item($row, 7).Value = 80
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
